$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: paragraph Range excluding the trailing paragraph-mark character, so
# that font / text edits don't bleed into the next paragraph.
# ---------------------------------------------------------------------------
function Get-TextRange($para) {
    $full = $para.Range
    return $d.Range($full.Start, $full.End - 1)
}

# ---------------------------------------------------------------------------
# 1) Paragraph "3)"  -> apply Times New Roman
# ---------------------------------------------------------------------------
$p16 = $d.Paragraphs(16)
$p16.Range.Font.Name = "Times New Roman"
$p16.Range.Font.NameBi = "Times New Roman"

# ---------------------------------------------------------------------------
# 2) Paragraph "File Name: queries__allAgg.js" -> apply Times New Roman
# ---------------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
$p17.Range.Font.Name = "Times New Roman"
$p17.Range.Font.NameBi = "Times New Roman"

# ---------------------------------------------------------------------------
# 3) Paragraph "With Index: Avg time - 1272.139ms"
#    -> replace text with "Without Index avg run time: 25865.3529"
# ---------------------------------------------------------------------------
$p18 = $d.Paragraphs(18)
$t18 = Get-TextRange $p18
$t18.Text = "Without Index avg run time: 25865.3529"

$p18 = $d.Paragraphs(18)
$p18.Range.Font.Name = "Times New Roman"
$p18.Range.Font.NameBi = "Times New Roman"
$p18.Range.Font.NameFarEast = "Times New Roman"

$t18 = Get-TextRange $p18
$t18.Font.Color = 2236962
$t18.Shading.BackgroundPatternColor = 16777215
$t18.Shading.ForegroundPatternColor = -16777216
$t18.Shading.Texture = 0

# ---------------------------------------------------------------------------
# 4) Paragraph "Without Index: Avg time - 1404.174ms"
#    -> replace text with "After Index on email ids: 21755.6325"
# ---------------------------------------------------------------------------
$p19 = $d.Paragraphs(19)
$t19 = Get-TextRange $p19
$t19.Text = "After Index on email ids: 21755.6325"

$p19 = $d.Paragraphs(19)
$p19.Range.Font.Name = "Times New Roman"
$p19.Range.Font.NameBi = "Times New Roman"
$p19.Range.Font.NameFarEast = "Times New Roman"

$t19 = Get-TextRange $p19
$t19.Font.Color = 2236962
$t19.Shading.BackgroundPatternColor = 16777215
$t19.Shading.ForegroundPatternColor = -16777216
$t19.Shading.Texture = 0

# ---------------------------------------------------------------------------
# 5) Paragraph "100 Aggregate Query" -> apply Times New Roman
# ---------------------------------------------------------------------------
$p20 = $d.Paragraphs(20)
$p20.Range.Font.Name = "Times New Roman"
$p20.Range.Font.NameBi = "Times New Roman"

# ---------------------------------------------------------------------------
# 6) Move the "_GoBack" bookmark from the end of the document to the end of
#    the (now-updated) "After Index on email ids" paragraph.
# ---------------------------------------------------------------------------
$bms = $d.Bookmarks
if ($bms.Exists("_GoBack")) {
    $old = $bms.Item("_GoBack")
    $old.Delete()
}

$p19 = $d.Paragraphs(19)
$target = Get-TextRange $p19
$bms.Add("_GoBack", $target)

Write-Host "done"
